$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.802.44'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.878.20'
$ws.Range('E3').Value = '  -0.81%  '
$ws.Range('D4').Value = '''1.004'
$ws.Range('E4').Value = '  -0.30%  '
$ws.Range('D5').Value = '''325.01'
$ws.Range('E5').Value = '  -0.62%  '
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Value = '''0.4595'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '''0.3880'
$ws.Range('E8').Value = '  +0.37%  '
$ws.Range('D9').Value = '''0.07864'
$ws.Range('E9').Value = '  -0.12%  '
$ws.Range('D10').Value = '''0.9839'
$ws.Range('E10').Value = '  -1.77%  '
$ws.Range('D11').Value = '''21.74'
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('D12').Value = '1.865.22'
$ws.Range('E12').Value = '  -1.16%  '
$ws.Range('D13').Value = '''7.005'
$ws.Range('E13').Value = '  -0.40%  '
$ws.Range('D14').Value = '''5.662'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').Value = '''0.06936'
$ws.Range('E15').Value = '  -0.26%  '
$ws.Range('D16').Value = '''88.29'
$ws.Range('E16').Value = '  +1.23%  '
$ws.Range('E17').Value = '  -0.34%  '
$ws.Range('D18').Value = '''0.000009966'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '''16.96'
$ws.Range('E19').Value = '  -0.26%  '
$ws.Range('D21').Value = '28.837.14'
$ws.Range('E21').Value = '  +0.48%  '
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('E23').Value = '  -0.46%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''2.090'
$ws.Range('E24').Value = '  +1.33%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''155.43'
$ws.Range('E25').Value = '  +0.48%  '
$ws.Range('B26').Value = 'EthereumClassic'
$ws.Range('C26').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D26').Value = '''19.28'
$ws.Range('E26').Value = '  -0.07%  '
$ws.Range('B27').Value = 'InternetComputer(DFINITY)'
$ws.Range('C27').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D27').Value = '''5.958'
$ws.Range('E27').Value = '  +2.66%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '''1.930'
$ws.Range('E28').Value = '  +1.25%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D29').Value = '''117.40'
$ws.Range('E29').Value = '  -0.74%  '
$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').Value = '''0.09332'
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''0.9027'
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '''5.268'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''1.326'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.266'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D35').Value = '''1.191'
$ws.Range('E35').Value = '  +2.75%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.05758'
$ws.Range('E36').Value = '  +0.47%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02071'
$ws.Range('E37').Value = '  +0.53%  '
$ws.Range('B38').Value = 'Frax'
$ws.Range('C38').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D38').Value = '''1.003'
$ws.Range('E38').Value = '  -0.21%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''7.667'
$ws.Range('E39').Value = '  -0.61%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5654'
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = '''0.1767'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').Value = '''9.681'
$ws.Range('E42').Value = '  -0.49%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '''2.261'
$ws.Range('E43').Value = '  +3.16%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '''11.86'
$ws.Range('E44').Value = '  +1.47%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').Value = '''0.5345'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '''0.07037'
$ws.Range('E46').Value = '  -1.85%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '''1.848'
$ws.Range('E47').Value = '  +1.17%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '''112.82'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').Value = '''2.519'
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D50').Value = '''1.064'
$ws.Range('E50').Value = '  -5.10%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').Value = '''70.70'
$ws.Range('E51').Value = '  +0.63%  '
